$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 677; existing rows 677-723 shift down to 679-725.
$ws.Rows("677:678").Insert()

# New row 677 (week of 2022-02-18 / serial 44610)
$ws.Range("A677").Value2 = 6
$ws.Range("B677").Value2 = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range("C677").Value2 = 'Metropolitana'
$ws.Range("D677").Value2 = 44610
$ws.Range("E677").Value2 = 13
$ws.Range("F677").Value2 = 100112040
$ws.Range("G677").Value2 = 'Cilantro'
$ws.Range("H677").Value2 = 'Sin especificar'
$ws.Range("I677").Value2 = 'Primera'
$ws.Range("J677").Value2 = 620
$ws.Range("K677").Value2 = 8000
$ws.Range("L677").Value2 = 9000
$ws.Range("M677").Value2 = 8403
$ws.Range("N677").Value2 = '$/caja 36 atados'
$ws.Range("O677").Value2 = 'Región Metropolitana'
$ws.Range("P677").Value2 = 233
$ws.Range("Q677").Value2 = 36
$ws.Range("R677").Value2 = 'Hortaliza'

# New row 678 (week of 2022-02-18 / serial 44610)
$ws.Range("A678").Value2 = 6
$ws.Range("B678").Value2 = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range("C678").Value2 = 'Metropolitana'
$ws.Range("D678").Value2 = 44610
$ws.Range("E678").Value2 = 13
$ws.Range("F678").Value2 = 100112040
$ws.Range("G678").Value2 = 'Cilantro'
$ws.Range("H678").Value2 = 'Sin especificar'
$ws.Range("I678").Value2 = 'Primera'
$ws.Range("J678").Value2 = 340
$ws.Range("K678").Value2 = 15000
$ws.Range("L678").Value2 = 16000
$ws.Range("M678").Value2 = 15441
$ws.Range("N678").Value2 = '$/docena de atados'
$ws.Range("O678").Value2 = 'Región Metropolitana'
$ws.Range("P678").Value2 = 5147
$ws.Range("Q678").Value2 = 3
$ws.Range("R678").Value2 = 'Hortaliza'
